# Weekly update: insert a new Coliflor / Macroferia Regional de Talca
# price record as row 138, shifting the existing rows (old 138..161) down
# by one (new rows 139..162).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 138 (pushes rows 138..161 down to 139..162,
# inheriting formatting from the row above, same as Excel's native
# "Insert Row" command).
$ws.Rows.Item(138).Insert()

# Populate the newly inserted row 138 with the new weekly record.
$ws.Cells.Item(138, 1).Value2 = 5
$ws.Cells.Item(138, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(138, 3).Value2 = "Maule"
$ws.Cells.Item(138, 4).Value2 = 44505
$ws.Cells.Item(138, 5).Value2 = 7
$ws.Cells.Item(138, 6).Value2 = 100112008
$ws.Cells.Item(138, 7).Value2 = "Coliflor"
$ws.Cells.Item(138, 8).Value2 = "Sin especificar"
$ws.Cells.Item(138, 9).Value2 = "Primera"
$ws.Cells.Item(138, 10).Value2 = 6000
$ws.Cells.Item(138, 11).Value2 = 500
$ws.Cells.Item(138, 12).Value2 = 500
$ws.Cells.Item(138, 13).Value2 = 500
$ws.Cells.Item(138, 14).Value2 = "$/unidad"
$ws.Cells.Item(138, 15).Value2 = "Región del Maule"
$ws.Cells.Item(138, 16).Value2 = 500
$ws.Cells.Item(138, 17).Value2 = 1
$ws.Cells.Item(138, 18).Value2 = "Hortaliza"
